# ER working within FF for multi-cycle simulations
#
# - Row 7: add "Problem Description" (C7) and "Fix Date" (D7)
# - Row 8 (new): Issue Date (A8), Problem Description (B8), Fix (C8), Relevant Code (E8)
# - Move the active selection to C10 (and drop the old topLeftCell/selection state)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 additions ---------------------------------------------------
$ws.Range("C7").Value = "Conversion of initial temperatures to Kelvin was done twice"

# Reuse the existing date formatting (style) from another date cell so the
# new date cell matches the workbook's established "Fix Date" / "Issue Date"
# number format instead of creating a brand new one.
[void]$ws.Range("D2").Copy()
[void]$ws.Range("D7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D7").Value = 44803              # 8/30/2022

# --- Row 8 (brand new row) ----------------------------------------------
[void]$ws.Range("A2").Copy()
[void]$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A8").Value = 44803              # 8/30/2022

$ws.Range("B8").Value = "Simulation stalls when performing running sequential cycles"
$ws.Range("C8").Value = "Seems to be a meshing issue… getting stalled in expm.m"
$ws.Range("E8").Value = "FF.m"

$excel.CutCopyMode = $false

# --- Selection / view state ---------------------------------------------
[void]$ws.Range("C10").Select()
